$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Bump the published term version and its date (per commit: "Added 1.1.0 of term")
$ws.Range("B3").Value = "1.1.0"
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"

# Re-apply wrap text so the cellXfs entries for these styled ranges carry
# applyAlignment="true" alongside the existing vertical="top" wrapText="true"
$ws.Range("A1:B14").WrapText = $true
